$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("AR4").Value = 1.8
$ws.Range("AS4").Value = 2.05
$ws.Range("M4").Value = 1.08
$ws.Range("N4").Value = 8
$ws.Range("S4").Value = 3.9
$ws.Range("T4").Value = 1.26
$ws.Range("AB5").Value = 7.5
$ws.Range("AD5").Value = 17
$ws.Range("AL5").Value = 8
$ws.Range("AQ5").Value = 67
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4.5
$ws.Range("AA6").Value = 5
$ws.Range("AB6").Value = 8
$ws.Range("AD6").Value = 17
$ws.Range("AH6").Value = 6.5
$ws.Range("AL6").Value = 8
$ws.Range("AM6").Value = 19
$ws.Range("G6").Value = 2.05
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 2.88
$ws.Range("L6").Value = 5
$ws.Range("Y6").Value = 2.38
$ws.Range("Z6").Value = 1.53
$ws.Range("AB7").Value = 29
$ws.Range("AD7").Value = 67
$ws.Range("AF7").Value = 51
$ws.Range("AH7").Value = 7
$ws.Range("AJ7").Value = 67
$ws.Range("AM7").Value = 6.5
$ws.Range("AN7").Value = 8.5
$ws.Range("AO7").Value = 11
$ws.Range("G7").Value = 6.5
$ws.Range("H7").Value = 3.7
$ws.Range("I7").Value = 1.57
$ws.Range("J7").Value = 6.5
$ws.Range("L7").Value = 2.2
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 8
$ws.Range("U7").Value = 3.75
$ws.Range("V7").Value = 1.25
$ws.Range("Y7").Value = 2.1
$ws.Range("Z7").Value = 1.67
$ws.Range("AR9").Value = 1.71
$ws.Range("AS9").Value = 2.1
$ws.Range("N9").Value = 8
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.6
$ws.Range("S9").Value = 3.5
$ws.Range("T9").Value = 1.29
$ws.Range("AA10").Value = 7.5
$ws.Range("AC10").Value = 13
$ws.Range("AD10").Value = 34
$ws.Range("AE10").Value = 34
$ws.Range("AG10").Value = 6
$ws.Range("AL10").Value = 6
$ws.Range("AM10").Value = 10
$ws.Range("AO10").Value = 23
$ws.Range("AP10").Value = 23
$ws.Range("G10").Value = 3.2
$ws.Range("I10").Value = 2.45
$ws.Range("J10").Value = 4
$ws.Range("L10").Value = 3.25
$ws.Range("Y10").Value = 2.2
$ws.Range("Z10").Value = 1.62
$ws.Range("AA11").Value = 5.6
$ws.Range("AC11").Value = 13
$ws.Range("AD11").Value = 45
$ws.Range("AE11").Value = 45
$ws.Range("AG11").Value = 4
$ws.Range("AJ11").Value = 200
$ws.Range("H11").Value = 2.45
$ws.Range("J11").Value = 3.9
$ws.Range("N11").Value = 4
$ws.Range("O11").Value = 1.8
$ws.Range("P11").Value = 1.91
$ws.Range("Q11").Value = 3.35
$ws.Range("R11").Value = 1.28
$ws.Range("U11").Value = 6.2
$ws.Range("W11").Value = 1.78
$ws.Range("X11").Value = 1.93
$ws.Range("AA12").Value = 6.7
$ws.Range("AB12").Value = 15
$ws.Range("AC12").Value = 12
$ws.Range("AE12").Value = 40
$ws.Range("AF12").Value = 60
$ws.Range("AG12").Value = 4.65
$ws.Range("AI12").Value = 18
$ws.Range("AJ12").Value = 120
$ws.Range("AL12").Value = 5.9
$ws.Range("AM12").Value = 11.5
$ws.Range("AN12").Value = 10.5
$ws.Range("AO12").Value = 32
$ws.Range("AP12").Value = 29
$ws.Range("AQ12").Value = 50
$ws.Range("H12").Value = 2.57
$ws.Range("I12").Value = 2.62
$ws.Range("J12").Value = 3.95
$ws.Range("K12").Value = 1.8
$ws.Range("M12").Value = 1.16
$ws.Range("N12").Value = 4.65
$ws.Range("O12").Value = 1.62
$ws.Range("P12").Value = 2.18
$ws.Range("Q12").Value = 2.8
$ws.Range("R12").Value = 1.38
$ws.Range("U12").Value = 5.1
$ws.Range("V12").Value = 1.13
$ws.Range("W12").Value = 1.62
$ws.Range("X12").Value = 2.15
$ws.Range("Y12").Value = 2.15
$ws.Range("Z12").Value = 1.62
$ws.Range("AA13").Value = 6.2
$ws.Range("AB13").Value = 11.75
$ws.Range("AC13").Value = 10.25
$ws.Range("AD13").Value = 32
$ws.Range("AE13").Value = 27
$ws.Range("AG13").Value = 4.9
$ws.Range("AH13").Value = 5.3
$ws.Range("AI13").Value = 17
$ws.Range("AJ13").Value = 110
$ws.Range("AL13").Value = 6.7
$ws.Range("AM13").Value = 14.5
$ws.Range("AN13").Value = 11.75
$ws.Range("AO13").Value = 45
$ws.Range("AP13").Value = 37
$ws.Range("AQ13").Value = 55
$ws.Range("G13").Value = 2.62
$ws.Range("H13").Value = 2.62
$ws.Range("I13").Value = 3.15
$ws.Range("J13").Value = 3.3
$ws.Range("K13").Value = 1.83
$ws.Range("L13").Value = 3.85
$ws.Range("N13").Value = 4.9
$ws.Range("P13").Value = 2.25
$ws.Range("Q13").Value = 2.67
$ws.Range("U13").Value = 4.8
$ws.Range("AA14").Value = 5.3
$ws.Range("AB14").Value = 5.5
$ws.Range("AD14").Value = 8.25
$ws.Range("AE14").Value = 13
$ws.Range("AG14").Value = 6.9
$ws.Range("AH14").Value = 8.5
$ws.Range("AI14").Value = 26
$ws.Range("AL14").Value = 17.5
$ws.Range("AM14").Value = 60
$ws.Range("AN14").Value = 28
$ws.Range("AO14").Value = 250
$ws.Range("AP14").Value = 120
$ws.Range("G14").Value = 1.38
$ws.Range("H14").Value = 4.15
$ws.Range("I14").Value = 8.75
$ws.Range("J14").Value = 1.87
$ws.Range("K14").Value = 2.25
$ws.Range("L14").Value = 7.8
$ws.Range("M14").Value = 1.07
$ws.Range("N14").Value = 6.9
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.05
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 1.72
$ws.Range("U14").Value = 3.35
$ws.Range("V14").Value = 1.28
$ws.Range("W14").Value = 1.4
$ws.Range("X14").Value = 2.72
$ws.Range("Y14").Value = 2.32
$ws.Range("Z14").Value = 1.55
$ws.Range("AC16").Value = 9
$ws.Range("AD16").Value = 12
$ws.Range("AE16").Value = 17
$ws.Range("AF16").Value = 41
$ws.Range("AG16").Value = 7.5
$ws.Range("AH16").Value = 7
$ws.Range("AL16").Value = 11
$ws.Range("AM16").Value = 26
$ws.Range("G16").Value = 1.67
$ws.Range("H16").Value = 3.5
$ws.Range("I16").Value = 5.75
$ws.Range("J16").Value = 2.38
$ws.Range("K16").Value = 2.05
$ws.Range("L16").Value = 6
$ws.Range("M16").Value = 1.08
$ws.Range("N16").Value = 8
$ws.Range("O16").Value = 1.44
$ws.Range("P16").Value = 2.75
$ws.Range("Q16").Value = 2.38
$ws.Range("R16").Value = 1.57
$ws.Range("U16").Value = 4.5
$ws.Range("V16").Value = 1.2
$ws.Range("Y16").Value = 2.25
$ws.Range("Z16").Value = 1.57
$ws.Range("AA17").Value = 7.5
$ws.Range("AL17").Value = 6
$ws.Range("AM17").Value = 10
$ws.Range("AN17").Value = 10
$ws.Range("AO17").Value = 23
$ws.Range("G17").Value = 3.1
$ws.Range("I17").Value = 2.4
$ws.Range("J17").Value = 4
$ws.Range("L17").Value = 3.25
$ws.Range("O17").Value = 1.53
$ws.Range("P17").Value = 2.5
$ws.Range("AK18").Value = 451
$ws.Range("H18").Value = 3.75
$ws.Range("AB20").Value = 10
$ws.Range("AC20").Value = 9
$ws.Range("AM20").Value = 17
$ws.Range("AN20").Value = 12
$ws.Range("G20").Value = 2.15
$ws.Range("H20").Value = 3.4
$ws.Range("I20").Value = 3.3
$ws.Range("J20").Value = 2.88
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = 1.06
$ws.Range("N20").Value = 10
$ws.Range("Q20").Value = 2.05
$ws.Range("R20").Value = 1.8
$ws.Range("U20").Value = 3.5
$ws.Range("V20").Value = 1.3
$ws.Range("AB21").Value = 9
$ws.Range("AD21").Value = 17
$ws.Range("AE21").Value = 17
$ws.Range("AN21").Value = 13
$ws.Range("G21").Value = 1.91
$ws.Range("I21").Value = 4.1
$ws.Range("Q21").Value = 2.1
$ws.Range("R21").Value = 1.73
$ws.Range("AB23").Value = 8
$ws.Range("AD23").Value = 17
$ws.Range("AH23").Value = 6
$ws.Range("G23").Value = 1.95
$ws.Range("H23").Value = 3.1
$ws.Range("I23").Value = 4.33
$ws.Range("K23").Value = 1.91
$ws.Range("L23").Value = 5
$ws.Range("AA24").Value = 7
$ws.Range("AB24").Value = 8
$ws.Range("AF24").Value = 26
$ws.Range("AL24").Value = 12
$ws.Range("AM24").Value = 23
$ws.Range("AN24").Value = 15
$ws.Range("G24").Value = 1.75
$ws.Range("H24").Value = 3.6
$ws.Range("I24").Value = 4.5
$ws.Range("J24").Value = 2.4
$ws.Range("L24").Value = 5
$ws.Range("Q24").Value = 2
$ws.Range("R24").Value = 1.8
$ws.Range("U24").Value = 3.5
$ws.Range("V24").Value = 1.29
$ws.Range("Y24").Value = 1.83
$ws.Range("Z24").Value = 1.83
